# Adds a new "NOF" (Number Of Functions) sheet, mirroring the layout of the
# existing "NOS" sheet, with its own clustered-column chart; also removes the
# title from the NOS chart (it becomes just "Number of Statements" -> no
# title, matching the new untitled NOF chart) and updates the remembered
# selection/active-sheet state left over from the editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Strip the title off the pre-existing NOS bar chart.
# ---------------------------------------------------------------------------
$wsNOS = $wb.Worksheets.Item("NOS")
$coNOS = $wsNOS.ChartObjects().Item(1)
$coNOS.Chart.HasTitle = $false

# ---------------------------------------------------------------------------
# 2. Refresh the cursor / selection left on the other sheets.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Average_cyclomatic")
$ws1.Range("A1:K3").Select()

$ws2 = $wb.Worksheets.Item("Comment_ratio")
$ws2.Range("K6").Select()

$ws3 = $wb.Worksheets.Item("NOS")
$ws3.Range("B11").Select()

# ---------------------------------------------------------------------------
# 3. Add the new "NOF" worksheet after NOS, populate its data table.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "NOF"

$projects = @("Connect-four", "Gomoku", "Tetris", "Mario-dodger", "Video Player", "Audio Player", "Video to MP3 Converter", "Shareit", "Messenger", "Video downloader")
$cols = @("C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
for ($i = 0; $i -lt $projects.Length; $i++) {
    $ws4.Range($cols[$i] + "2").Value = $projects[$i]
}

$ws4.Range("B3").Value = "JavaScript"
$jsValues = @(19, 64, 204, 1124, 98, 214, 8, 429, 429, 10)
for ($i = 0; $i -lt $jsValues.Length; $i++) {
    $ws4.Range($cols[$i] + "3").Value = $jsValues[$i]
}

$ws4.Range("B4").Value = "TypeScript"
$tsValues = @(37, 113, 7, 28, 232, 232, 20, 162, 185, 145)
for ($i = 0; $i -lt $tsValues.Length; $i++) {
    $ws4.Range($cols[$i] + "4").Value = $tsValues[$i]
}

# Match the bestFit column widths used on the sibling "NOS" sheet (shifted
# one column to the right: NOS!A:K -> NOF!B:L).
$nosWidths = @{}
foreach ($col in @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")) {
    $nosWidths[$col] = $ws3.Columns($col).ColumnWidth
}
$colMap = @{"B" = "A"; "C" = "B"; "D" = "C"; "E" = "D"; "F" = "E"; "G" = "F"; "H" = "G"; "I" = "H"; "J" = "I"; "K" = "J"; "L" = "K"}
foreach ($nofCol in $colMap.Keys) {
    $ws4.Columns($nofCol).ColumnWidth = $nosWidths[$colMap[$nofCol]]
}

# ---------------------------------------------------------------------------
# 4. Build a clustered-column chart for NOF, mirroring the NOS chart (two
#    series - JavaScript / TypeScript - across the ten sample projects) but
#    without a title, same as the now-untitled NOS chart.
# ---------------------------------------------------------------------------
$co = $ws4.ChartObjects().Add(169, 98, 467, 274)
$chart = $co.Chart
$chart.ChartType = 51

$chart.SeriesCollection().NewSeries()
$sJs = $chart.SeriesCollection().Item(1)
$sJs.Name = "=NOF!`$B`$3"
$sJs.Values = "=NOF!`$C`$3:`$L`$3"
$sJs.XValues = "=NOF!`$C`$2:`$L`$2"

$chart.SeriesCollection().NewSeries()
$sTs = $chart.SeriesCollection().Item(2)
$sTs.Name = "=NOF!`$B`$4"
$sTs.Values = "=NOF!`$C`$4:`$L`$4"
$sTs.XValues = "=NOF!`$C`$2:`$L`$2"

$chart.HasTitle = $false
$chart.HasLegend = $true
$chart.Legend.Position = -4107

# ---------------------------------------------------------------------------
# 5. Leave the selection/active sheet on NOF, matching the saved session.
# ---------------------------------------------------------------------------
$ws4.Range("B2:L4").Select()
